$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: new entry (SINO 4, June/08/2021, Practical, Linked List topic)
$ws.Range("B13").Value = "June/08/2021"
$ws.Range("A13").Value = 4
$ws.Range("C13").Value = "Practical"
$ws.Range("D13").Value = "Linked List"
$ws.Range("E13").Value = "1.Delete all occurance of the number inputed"
$ws.Range("F13").Value = "Completed"
$ws.Range("H13").Value = "Completed"

# Row 14
$ws.Range("E14").Value = "2.Remove duplicate element from the list"
$ws.Range("F14").Value = "Completed"
$ws.Range("H14").Value = "Completed"

# Row 15
$ws.Range("E15").Value = "3.Delete first node of the list"
$ws.Range("F15").Value = "Completed"
$ws.Range("H15").Value = "Completed"

# Row 16
$ws.Range("E16").Value = "4.Find the middle element using two pointer algorithm"
$ws.Range("F16").Value = "Completed"

# Row 17
$ws.Range("E17").Value = "5.Find the loop in the linked list using two pointer algorithm"
$ws.Range("F17").Value = "Completed"

# Row 7: replace the old description text, and mark Sur_status Completed
$ws.Range("E7").Value = "1.Delete the element at the last"
$ws.Range("H7").Value = "Completed"
$ws.Range("H9").Value = "Completed"

# Resize Description column to fit the new, longer text
$ws.Columns.Item(5).AutoFit()

# Leave selection where the user last clicked while reviewing the new rows
$ws.Range("L14").Select()
